$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-7), keep header row 1 intact;
# the new dataset (rows 2-13) is written fresh below.
$ws.Range("A2:Q7").ClearContents()

# Row 2: TRIAL SERVICES COMPANY  LTD
$ws.Range("A2").Value = "3103498400000V"
$ws.Range("B2").Value = "TRIAL SERVICES COMPANY  LTD"
$ws.Range("C2").Value = "TIN"
$ws.Range("D2").Value = "TRIAL SERVICES COMPANY  LTD"
$ws.Range("E2").Value = "Nyabugogo :"
$ws.Range("F2").Value = "Nyabugogo :"
$ws.Range("G2").Value = "Nyabugogo :"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "1245 Kigali"
$ws.Range("I2").Value = "Rwanda"
$ws.Range("J2").Value = "trialservicesltd@yahoo.com"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "788382272"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "1197080004983020"
$ws.Range("O2").Value = "RWANDA"
$ws.Range("P2").Value = "COGEBANQUE"
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "102317498"

# Row 3: GREAT HOTEL KIYOVU Ltd
$ws.Range("A3").Value = "3117513300000X"
$ws.Range("B3").Value = "GREAT HOTEL KIYOVU Ltd"
$ws.Range("C3").Value = "TIN"
$ws.Range("D3").Value = "GREAT HOTEL KIYOVU Ltd"
$ws.Range("E3").Value = "Nyarugenge-Kigali"
$ws.Range("F3").Value = "Nyarugenge-Kigali"
$ws.Range("G3").Value = "Nyarugenge-Kigali"
$ws.Range("I3").Value = "Rwanda"
$ws.Range("J3").Value = "greathotel2050@gmail.com"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "788382272"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "1197080004983020"
$ws.Range("O3").Value = "RWANDA"
$ws.Range("P3").Value = "BK"
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "110329376"

# Row 4: ITAU AUDITORS AND DPA  LTD
$ws.Range("A4").Value = "3101669500000T"
$ws.Range("B4").Value = "ITAU AUDITORS AND DPA  LTD"
$ws.Range("C4").Value = "TIN"
$ws.Range("D4").Value = "ITAU AUDITORS AND DPA  LTD"
$ws.Range("E4").Value = "NYARUGENGE"
$ws.Range("F4").Value = "NYARUGENGE"
$ws.Range("G4").Value = "NYARUGENGE"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "4385"
$ws.Range("I4").Value = "Rwanda"
$ws.Range("J4").Value = "itau_dpa@yahoo.fr"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "0788307360"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = "1196180001858031"
$ws.Range("O4").Value = "RWANDA"
$ws.Range("P4").Value = "EQUITY BANK"
$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "101318041"

# Row 5: ITAU AUDITORS  LTD
$ws.Range("A5").Value = "3102294700000A"
$ws.Range("B5").Value = "ITAU AUDITORS  LTD"
$ws.Range("C5").Value = "TIN"
$ws.Range("D5").Value = "ITAU AUDITORS  LTD"
$ws.Range("E5").Value = "NYARUGENGE"
$ws.Range("F5").Value = "NYARUGENGE"
$ws.Range("G5").Value = "NYARUGENGE"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "4385"
$ws.Range("I5").Value = "Rwanda"
$ws.Range("J5").Value = "itau_dpa@yahoo.fr"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "0788677410"
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "KE021952"
$ws.Range("O5").Value = "KENYA"
$ws.Range("P5").Value = "I&M Bank"
$ws.Range("Q5").NumberFormat = "@"
$ws.Range("Q5").Value = "102007253"

# Row 6: DIRECT SERVICES  LTD
$ws.Range("A6").Value = "3103164500000M"
$ws.Range("B6").Value = "DIRECT SERVICES  LTD"
$ws.Range("C6").Value = "TIN"
$ws.Range("D6").Value = "DIRECT SERVICES  LTD"
$ws.Range("E6").Value = "Karubanda :"
$ws.Range("F6").Value = "Karubanda :"
$ws.Range("G6").Value = "Karubanda :"
$ws.Range("I6").Value = "Rwanda"
$ws.Range("J6").Value = "rupaccy@gmail.com"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "0788359736"
$ws.Range("N6").NumberFormat = "@"
$ws.Range("N6").Value = "1198980180555085"
$ws.Range("O6").Value = "RWANDA"
$ws.Range("P6").Value = "COGEBANQUE"
$ws.Range("Q6").NumberFormat = "@"
$ws.Range("Q6").Value = "102190716"

# Row 7: RUPA BUSINESS COMPANY Ltd
$ws.Range("A7").Value = "3117465400000Y"
$ws.Range("B7").Value = "RUPA BUSINESS COMPANY Ltd"
$ws.Range("C7").Value = "TIN"
$ws.Range("D7").Value = "RUPA BUSINESS COMPANY Ltd"
$ws.Range("E7").Value = "Bugesera, Nyamata"
$ws.Range("F7").Value = "Bugesera, Nyamata"
$ws.Range("G7").Value = "Bugesera, Nyamata"
$ws.Range("I7").Value = "Rwanda"
$ws.Range("J7").Value = "rupaccy@gmail.com"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "788359736"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "1198980180555085"
$ws.Range("O7").Value = "RWANDA"
$ws.Range("P7").Value = "COGEBANQUE"
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "110179209"

# Row 8: ULTIMATE DEVELOPERS Ltd
$ws.Range("A8").Value = "3104214900000Z"
$ws.Range("B8").Value = "ULTIMATE DEVELOPERS Ltd"
$ws.Range("C8").Value = "TIN"
$ws.Range("D8").Value = "ULTIMATE DEVELOPERS Ltd"
$ws.Range("E8").Value = "Nyarugenge: RSSB Building-Tower II-9th Floor"
$ws.Range("F8").Value = "Nyarugenge: RSSB Building-Tower II-9th Floor"
$ws.Range("G8").Value = "Nyarugenge: RSSB Building-Tower II-9th Floor"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "5516"
$ws.Range("I8").Value = "Rwanda"
$ws.Range("J8").Value = "d.murwanashyaka@udl.rw"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "0780867765"
$ws.Range("M8").Value = "www.udl.rw"
$ws.Range("N8").NumberFormat = "@"
$ws.Range("N8").Value = "1198980178061009"
$ws.Range("O8").Value = "RWANDA"
$ws.Range("P8").Value = "BK"
$ws.Range("Q8").NumberFormat = "@"
$ws.Range("Q8").Value = "102514147"

# Row 9: INUMA TECHNOLOGY Ltd
$ws.Range("A9").Value = "3114859500000P"
$ws.Range("B9").Value = "INUMA TECHNOLOGY Ltd"
$ws.Range("C9").Value = "TIN"
$ws.Range("D9").Value = "INUMA TECHNOLOGY Ltd"
$ws.Range("E9").Value = "KN 2 ST Rubangura House(Underground Floor-Door 108)"
$ws.Range("F9").Value = "KN 2 ST Rubangura House(Underground Floor-Door 108)"
$ws.Range("G9").Value = "KN 2 ST Rubangura House(Underground Floor-Door 108)"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "+250"
$ws.Range("I9").Value = "Rwanda"
$ws.Range("J9").Value = "mugabejosue@gmail.com"
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = "0780867765"
$ws.Range("M9").Value = "www.inumatechnology.com"
$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "1199470011111103"
$ws.Range("O9").Value = "Rwanda"
$ws.Range("P9").Value = "Bank of Africa Rwanda"
$ws.Range("Q9").NumberFormat = "@"
$ws.Range("Q9").Value = "107530868"

# Row 10: KALISA STEPHEN
$ws.Range("B10").Value = "KALISA STEPHEN"
$ws.Range("C10").Value = "TIN"
$ws.Range("D10").Value = "KALISA STEPHEN"
$ws.Range("E10").Value = "KIGALI CITY"
$ws.Range("F10").Value = "KIGALI CITY"
$ws.Range("G10").Value = "KIGALI CITY"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "1975"
$ws.Range("I10").Value = "Rwanda"
$ws.Range("J10").Value = "stephen.cp12@gmail.com"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "0788306945"
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "1197880014571139"
$ws.Range("O10").Value = "Rwandan"
$ws.Range("P10").Value = "EQUITY BANK"
$ws.Range("Q10").NumberFormat = "@"
$ws.Range("Q10").Value = "103563793"

# Row 11: ROCK STARS CONSULTANCY & SUPPLIERS Ltd
$ws.Range("A11").Value = "3104239700000H"
$ws.Range("B11").Value = "ROCK STARS CONSULTANCY & SUPPLIERS Ltd"
$ws.Range("C11").Value = "TIN"
$ws.Range("D11").Value = "ROCK STARS CONSULTANCY & SUPPLIERS Ltd"
$ws.Range("E11").Value = "Kamashashi :"
$ws.Range("F11").Value = "Kamashashi :"
$ws.Range("G11").Value = "Kamashashi :"
$ws.Range("I11").Value = "Rwanda"
$ws.Range("J11").Value = "rockconsult13@yahoo.co.uk"
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = "0788306945"
$ws.Range("N11").NumberFormat = "@"
$ws.Range("N11").Value = "1197880014571139"
$ws.Range("O11").Value = "RWANDA"
$ws.Range("P11").Value = "BK"
$ws.Range("Q11").NumberFormat = "@"
$ws.Range("Q11").Value = "102520083"

# Row 12: IBIGABIRO HOTEL Ltd
$ws.Range("A12").Value = "3115399100000H"
$ws.Range("B12").Value = "IBIGABIRO HOTEL Ltd"
$ws.Range("C12").Value = "TIN"
$ws.Range("D12").Value = "IBIGABIRO HOTEL Ltd"
$ws.Range("E12").Value = "KAGANO/NYAMASHEKE"
$ws.Range("F12").Value = "KAGANO/NYAMASHEKE"
$ws.Range("G12").Value = "KAGANO/NYAMASHEKE"
$ws.Range("I12").Value = "Rwanda"
$ws.Range("J12").Value = "judith.uwankwera@yahoo.com"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = "788777282"
$ws.Range("N12").NumberFormat = "@"
$ws.Range("N12").Value = "1196970030416098"
$ws.Range("O12").Value = "RWANDA"
$ws.Range("P12").Value = "COGEBANQUE"
$ws.Range("Q12").NumberFormat = "@"
$ws.Range("Q12").Value = "107931976"

# Row 13: SPEC  LTD
$ws.Range("A13").Value = "3104017800000C"
$ws.Range("B13").Value = "SPEC  LTD"
$ws.Range("C13").Value = "TIN"
$ws.Range("D13").Value = "SPEC  LTD"
$ws.Range("E13").Value = "Nyamasheke"
$ws.Range("F13").Value = "Nyamasheke"
$ws.Range("G13").Value = "Nyamasheke"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "12 , Nyamasheke"
$ws.Range("I13").Value = "Rwanda"
$ws.Range("J13").Value = "judith.uwankwera@yahoo.com"
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = "788777282"
$ws.Range("N13").NumberFormat = "@"
$ws.Range("N13").Value = "1196970030416098"
$ws.Range("O13").Value = "RWANDA"
$ws.Range("P13").Value = "COGEBANQUE"
$ws.Range("Q13").NumberFormat = "@"
$ws.Range("Q13").Value = "102461662"
